$wb = $excel.ActiveWorkbook

# --- Update AddCustomerTest sheet: new rows for vicky/yadav/gaurav ---
$ac = $wb.Worksheets.Item("AddCustomerTest")
$ac.Range("A4").Value = "vicky"
$ac.Range("B5").Value = "yadav"
$ac.Range("A5").Value = "gaurav"
$ac.Range("B4").Value = "thopate"
$ac.Range("C4").Value = 412303
$ac.Range("C5").Value = 412301
$ac.Range("D4").Value = "Customer added successfully"
$ac.Range("D5").Value = "Customer added successfully"

# --- Update OpenAccountTest sheet: capitalize currency name ---
$oa = $wb.Worksheets.Item("OpenAccountTest")
$oa.Range("B2").Value = "Rupee"
$oa.Activate()
$oa.Range("B2").Select()

# --- Add the new test_suite sheet and move it to the front ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "test_suite"
$newSheet.Move($wb.Worksheets.Item(1))

$ts = $wb.Worksheets.Item("test_suite")
$ts.Range("A1").Value = "TCID"
$ts.Range("B1").Value = "Runmode"
$ts.Range("A2").Value = "BankManagerLoginTest"
$ts.Range("A3").Value = "AddCustomerTest"
$ts.Range("A4").Value = "OpenAccountTest"
$ts.Range("B2").Value = "Y"
$ts.Range("B3").Value = "Y"
$ts.Range("B4").Value = "N"

$ts.Columns(1).ColumnWidth = 24.33203125
$ts.Columns(2).ColumnWidth = 15.88671875

$ts.Activate()
$ts.Range("B1:B4").Select()

# --- Re-select AddCustomerTest as the active tab, matching activeTab index 1 ---
$ac = $wb.Worksheets.Item("AddCustomerTest")
$ac.Activate()
$ac.Range("E1:E7").Select()

Write-Host "done"
